$d = $word.ActiveDocument

# Table 3 ("Työmäärä toteutettu" = 40 for both Topi Uhtakari and Severi
# lillkåll) -> both become 38.
$t1 = $d.Tables.Item(3)
$t1.Rows.Item(2).Cells.Item(3).Range.Text = "38"
$t1.Rows.Item(3).Cells.Item(3).Range.Text = "38"

# Table 4 ("Työmäärä toteutettu" = 15 for Severi lillkåll row) -> becomes 25.
$t2 = $d.Tables.Item(4)
$t2.Rows.Item(3).Cells.Item(3).Range.Text = "25"

Write-Output "done"
